# Apply edits described by the commit: "xoa cac phan tong cua ti le chiet khau"
# (remove the discount-rate totals)
#
# Changes:
#  1. Sheet "Đơn sale chính": M4 (discount-rate total) 0.2 -> 0
#  2. Sheet "Đơn phụ phẫu 1": new data row inserted before the "Tổng" row,
#     and the "Tổng" row's totals updated accordingly.
#  3. Sheet "Lương": remove the "HỆ THỐNG" sub-section (and its grand-total
#     row), then update a handful of recalculated totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Đơn sale chính"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 13).Value = 0   # M4: 0.2 -> 0

# ---------------------------------------------------------------------
# Sheet 2: "Đơn phụ phẫu 1"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row at position 14, pushing the old "Tổng" row down to 15.
$ws2.Rows.Item(14).Insert()

$ws2.Cells.Item(14, 1).Value = "HD-LUXURY"
$ws2.Cells.Item(14, 2).Value = 590
# Leading apostrophe forces the date-like string to stay plain text instead
# of being auto-converted to a date serial number (matches the other rows,
# which store this column as plain text).
$ws2.Cells.Item(14, 3).Value = "'07-25-2024"
$ws2.Cells.Item(14, 4).Value = "CẦN THƠ"
$ws2.Cells.Item(14, 5).Value = "Trần Thị Lệ"
$ws2.Cells.Item(14, 6).Value = "Cá nhân"
$ws2.Cells.Item(14, 7).Value = "Nâng mũi"
$ws2.Cells.Item(14, 8).Value = "Lâm Hoàng Phú"
$ws2.Cells.Item(14, 9).Value = 100000

# Update the (now shifted) "Tổng" row 15 totals.
$ws2.Cells.Item(15, 2).Value = 13
$ws2.Cells.Item(15, 9).Value = 1050000

# ---------------------------------------------------------------------
# Sheet 3: "Lương"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Remove the "... tại HỆ THỐNG" detail block (rows 4-10).
$ws3.Range("A4:A10").EntireRow.Delete()

# Remove the "Tổng lương tại HỆ THỐNG" row (row 35 originally, now row 28
# after the previous 7-row deletion).
$ws3.Range("A28").EntireRow.Delete()

# Update the recalculated figures.
$ws3.Cells.Item(2, 2).Value = 24.5           # Ngày công
$ws3.Cells.Item(3, 2).Value = 857500         # Phụ cấp
$ws3.Cells.Item(4, 2).Value = 2625000        # Lương cơ bản tại CẦN THƠ
$ws3.Cells.Item(9, 2).Value = 1050000        # Công phụ phẫu 1 tại CẦN THƠ
$ws3.Cells.Item(28, 2).Value = 3582500       # Tổng lương tại CẦN THƠ
$ws3.Cells.Item(31, 2).Value = 3651300       # Tổng lương

# The row-delete above turns the originally-blank "Lương cơ bản tại ..."
# cells (rows 12 and 20) into literal 0s; restore them to blank so they
# stay empty as in the source data.
$ws3.Cells.Item(12, 2).ClearContents()       # Lương cơ bản tại LONG XUYÊN
$ws3.Cells.Item(20, 2).ClearContents()       # Lương cơ bản tại SÓC TRĂNG
